$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'66.137.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.90%  '

# Row 3
$ws.Range("D3").Value = "'2.987.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.54%  '

# Row 4
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
$ws.Range("D5").Value = "'577.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '

# Row 6
$ws.Range("D6").Value = "'161.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.94%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("D8").Value = "'0.516"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.20%  '

# Row 9
$ws.Range("D9").Value = "'2.983.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.56%  '

# Row 10
$ws.Range("D10").Value = "'6.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.98%  '

# Row 11
$ws.Range("D11").Value = "'0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.70%  '

# Row 12
$ws.Range("D12").Value = "'0.452"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.04%  '

# Row 13
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.30%  '

# Row 14
$ws.Range("D14").Value = "'34.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.55%  '

# Row 15
$ws.Range("E15").Value = '  -0.61%  '

# Row 16
$ws.Range("D16").Value = "'66.134.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.19%  '

# Row 17
$ws.Range("D17").Value = "'3.477.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.67%  '

# Row 18
$ws.Range("D18").Value = "'6.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.26%  '

# Row 19
$ws.Range("D19").Value = "'2.990.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.29%  '

# Row 20
$ws.Range("D20").Value = "'449.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.63%  '

# Row 21
$ws.Range("D21").Value = "'13.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.99%  '

# Row 22
$ws.Range("D22").Value = "'0.679"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.59%  '

# Row 23
$ws.Range("D23").Value = "'7.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.15%  '

# Row 24
$ws.Range("D24").Value = "'81.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.20%  '

# Row 25
$ws.Range("D25").Value = "'2.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.24%  '

# Row 26
$ws.Range("D26").Value = "'12.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.24%  '

# Row 27
$ws.Range("D27").Value = "'10.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.94%  '

# Row 28
$ws.Range("E28").Value = '  +0.03%  '

# Row 29
$ws.Range("D29").Value = "'8.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +10.06%  '

# Row 30
$ws.Range("E30").Value = '  +14.50%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = "'2.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.05%  '

# Row 32
$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").Value = "'0.0000102"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.52%  '

# Row 33
$ws.Range("D33").Value = "'26.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.93%  '

# Row 34
$ws.Range("D34").Value = "'0.109"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.37%  '

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.25%  '

# Row 36
$ws.Range("D36").Value = "'0.982"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.44%  '

# Row 37
$ws.Range("D37").Value = "'5.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.08%  '

# Row 38
$ws.Range("D38").Value = "'2.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.47%  '

# Row 39
$ws.Range("D39").Value = "'49.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.73%  '

# Row 40
$ws.Range("B40").Value = 'Arweave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D40").Value = "'43.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.98%  '

# Row 41
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").Value = "'2.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.10%  '

# Row 42
$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D42").Value = "'0.301"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.17%  '

# Row 43
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").Value = "'0.120"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.12%  '

# Row 44
$ws.Range("D44").Value = "'8.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.70%  '

# Row 45
$ws.Range("D45").Value = "'394.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.20%  '

# Row 46
$ws.Range("E46").Value = '  +3.51%  '

# Row 47
$ws.Range("D47").Value = "'2.741.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.07%  '

# Row 48
$ws.Range("D48").Value = "'132.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.56%  '

# Row 49
$ws.Range("E49").Value = '  +0.00%  '

# Row 50
$ws.Range("D50").Value = "'23.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.58%  '

# Row 51
$ws.Range("E51").Value = '  +2.32%  '
